$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) Status text: "Ready for handoff" -> "In Translation"
#    This string is shared by the Status cells on all three sheets:
#      Overview!E2  (zh-cn status)
#      Overview!F2  (de-de status)
#      zh-cn!C2     (Status column)
#      de-de!C2     (Status column)
# ---------------------------------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E2").Value = "In Translation"
$overview.Range("F2").Value = "In Translation"

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C2").Value = "In Translation"

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C2").Value = "In Translation"

# ---------------------------------------------------------------------------
# 2) Narrow the "zh-cn" / "de-de" status columns (report was regenerated with
#    a tighter auto-fit width):
#      Overview columns E & F
#      zh-cn column C
#      de-de column C
# ---------------------------------------------------------------------------
$overview.Range("E1:F1").EntireColumn.ColumnWidth = 12.5
$zhcn.Range("C1").EntireColumn.ColumnWidth = 12.5
$dede.Range("C1").EntireColumn.ColumnWidth = 12.5
